# Append two new result rows (12 and 13) that duplicate the last existing
# data row (row 11), only changing the Timestamp value in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template row with all the stats that should be repeated.
$template = $ws.Range("A11:BI11")

# New row 12: copy the template row, then overwrite the timestamp.
$row12 = $ws.Range("A12:BI12")
$template.Copy($row12)
$ws.Range("A12").Value = "2025-08-20 13:13:44"

# New row 13: copy the template row, then overwrite the timestamp.
$row13 = $ws.Range("A13:BI13")
$template.Copy($row13)
$ws.Range("A13").Value = "2025-08-20 13:17:56"
